$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "All" row (row 2) values with final prepared data
$ws.Range("B2").Value = 0.0125997860246464
$ws.Range("C2").Value = 0.00103818995711643
$ws.Range("D2").Value = 0.0241613820921764

# Add the new "along"/info_solidarity indicator column (E)
$ws.Range("E1").Value = "along"
$ws.Range("E2").Value = "info_solidarityTRUE"
$ws.Range("E3").Value = "info_solidarityTRUE"
$ws.Range("E4").Value = "info_solidarityTRUE"
$ws.Range("E5").Value = "info_solidarityTRUE"
$ws.Range("E6").Value = "info_solidarityTRUE"
$ws.Range("E7").Value = "info_solidarityTRUE"
$ws.Range("E8").Value = "info_solidarityTRUE"
$ws.Range("E9").Value = "info_solidarityTRUE"
$ws.Range("E10").Value = "info_solidarityTRUE"
$ws.Range("E11").Value = "info_solidarityTRUE"

# Insert a new "Russia" row before the current "Saudi Arabia" row (row 12),
# pushing Saudi Arabia and USA down by one row
$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = "Russia"
$ws.Range("B12").Value = -0.0068616082827698
$ws.Range("C12").Value = -0.0485694859957878
$ws.Range("D12").Value = 0.0348462694302482
$ws.Range("E12").Value = "info_solidarityTRUE"

# Fill in the new column for the rows that were pushed down
$ws.Range("E13").Value = "info_solidarityTRUE"
$ws.Range("E14").Value = "info_solidarityTRUE"
